$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '29.474.93'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = "'" + '1.875.34'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("D4").Value = "'" + '0.9989'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'" + '0.7152'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").Value = "'" + '241.50'
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("D7").Value = "'" + '0.9991'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = "'" + '0.07949'
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").Value = "'" + '0.3106'
$ws.Range("E9").Value = '  +3.01%  '
$ws.Range("E10").Value = '  +5.98%  '
$ws.Range("D11").Value = "'" + '0.08274'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = "'" + '0.7299'
$ws.Range("E12").Value = '  +3.46%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'" + '5.279'
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'" + '1.863.66'
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").Value = "'" + '91.16'
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").Value = "'" + '29.480.40'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = "'" + '5.928'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("D18").Value = "'" + '245.39'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = "'" + '13.34'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = "'" + '2.118.31'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").Value = "'" + '0.9993'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = "'" + '7.990'
$ws.Range("E23").Value = '  +6.76%  '
$ws.Range("D24").Value = "'" + '0.9986'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").Value = "'" + '0.1610'
$ws.Range("E25").Value = '  +13.70%  '
$ws.Range("D26").Value = "'" + '163.54'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = "'" + '9.046'
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("D28").Value = "'" + '18.30'
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("D29").Value = "'" + '1.358'
$ws.Range("E29").Value = '  -3.02%  '
$ws.Range("D30").Value = "'" + '1.489'
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("D31").Value = "'" + '4.387'
$ws.Range("E31").Value = '  +2.14%  '
$ws.Range("D33").Value = "'" + '0.05269'
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("D34").Value = "'" + '1.954'
$ws.Range("E34").Value = '  +2.19%  '
$ws.Range("D35").Value = "'" + '1.198'
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("D36").Value = "'" + '0.7276'
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("D37").Value = "'" + '2.673'
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = "'" + '0.01869'
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("D39").Value = "'" + '1.223.12'
$ws.Range("E39").Value = '  +6.40%  '
$ws.Range("D40").Value = "'" + '2.706'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = "'" + '0.9094'
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("D42").Value = "'" + '73.81'
$ws.Range("E42").Value = '  +5.56%  '
$ws.Range("D43").Value = "'" + '6.117'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").Value = "'" + '0.9991'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = "'" + '102.30'
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = "'" + '2.014.33'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'" + '0.5284'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E48").Value = '  +3.73%  '
$ws.Range("D49").Value = "'" + '2.919'
$ws.Range("E49").Value = '  +8.87%  '
$ws.Range("D50").Value = "'" + '0.00000000120'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").Value = "'" + '9.347'
$ws.Range("E51").Value = '  +2.16%  '
